$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.070284843444824
$ws.Range("B1").Value = 1.625010251998901
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.877082347869873
$ws.Range("E1").Value = 1.159711360931396
